$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "number_of_dwellings" baseline series in column B (rows 2-402, years 1600-2000)
# was recalibrated with a type-split multiplier that differs by period:
#   years 1600-1980 (rows   2-382): x1.4
#   years 1981-1990 (rows 383-392): x1.848
#   years 1991-2000 (rows 393-402): x1.54
# Rows 403-452 (years 2001-2050) are left untouched.

$ranges = @(
    @{ Range = "B2:B382";  Factor = 1.4   },
    @{ Range = "B383:B392"; Factor = 1.848 },
    @{ Range = "B393:B402"; Factor = 1.54  }
)

foreach ($item in $ranges) {
    $rng = $ws.Range($item.Range)
    $factor = $item.Factor
    foreach ($cell in $rng.Cells) {
        $cell.Value = $cell.Value2 * $factor
    }
}
